# FUNCTIONALITY: Wrote a new Suite.
# Add a new test-suite row (row 3) to the statistics sheet, with values
# mirroring the layout of the existing data row (row 2), and move the
# active selection to E4 (just below the newly written row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "Runtimes"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 9
$ws.Range("D3").Value = "Suited to Manual"
$ws.Range("E3").Value = "Timing"

$ws.Range("E4").Select() | Out-Null
